$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.771.43"
$ws.Range("E2").Value = "  -4.65%  "
$ws.Range("D3").Value = "2.310.52"
$ws.Range("E3").Value = "  -6.48%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.53"
$ws.Range("E5").Value = "  -4.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "83.99"
$ws.Range("E6").Value = "  -8.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.528"
$ws.Range("E7").Value = "  -3.87%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.480"
$ws.Range("E9").Value = "  -5.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0805"
$ws.Range("E10").Value = "  -5.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "29.72"
$ws.Range("E11").Value = "  -9.64%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "2.671.83"
$ws.Range("E13").Value = "  -6.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.36"
$ws.Range("E14").Value = "  -7.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.57"
$ws.Range("E15").Value = "  -5.91%  "
$ws.Range("D16").Value = "2.327.06"
$ws.Range("E16").Value = "  -5.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.747"
$ws.Range("E17").Value = "  -5.38%  "
$ws.Range("D18").Value = "39.768.23"
$ws.Range("E18").Value = "  -4.52%  "
$ws.Range("D19").Value = "0.0₃0896"
$ws.Range("E19").Value = "  -4.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.03"
$ws.Range("E20").Value = "  -6.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.09"
$ws.Range("E21").Value = "  -6.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.49"
$ws.Range("E22").Value = "  -6.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.59"
$ws.Range("E23").Value = "  -2.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.54"
$ws.Range("E24").Value = "  -7.77%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.79"
$ws.Range("E26").Value = "  -7.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.97"
$ws.Range("E27").Value = "  -7.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.20"
$ws.Range("E29").Value = "  -5.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.46"
$ws.Range("E30").Value = "  -4.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "150.91"
$ws.Range("E31").Value = "  -2.98%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.04"
$ws.Range("E33").Value = "  -7.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.45"
$ws.Range("E34").Value = "  -4.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0714"
$ws.Range("E35").Value = "  -6.47%  "
$ws.Range("E36").Value = "  -2.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0987"
$ws.Range("E37").Value = "  -4.18%  "
$ws.Range("E38").Value = "  -6.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.38"
$ws.Range("E39").Value = "  -9.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.69"
$ws.Range("E40").Value = "  -7.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.77"
$ws.Range("E41").Value = "  -5.62%  "
$ws.Range("E42").Value = "  -2.30%  "
$ws.Range("D43").Value = "1.935.91"
$ws.Range("E43").Value = "  -3.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0264"
$ws.Range("E44").Value = "  -6.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.31"
$ws.Range("E45").Value = "  -6.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.34"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.65"
$ws.Range("E47").Value = "  -10.36%  "
$ws.Range("D48").Value = "2.548.92"
$ws.Range("E48").Value = "  -6.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.01"
$ws.Range("E49").Value = "  -5.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.94"
$ws.Range("E50").Value = "  -7.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.77"
$ws.Range("E51").Value = "  -4.57%  "
